# Loan RBI, Variable Instalments
#
# On the "Repayment Schedule" sheet, a new blank column is inserted between
# the existing "In Advance" column (old column N) and the "Over Due" column
# (old column O). Everything from old column N onward shifts one column to
# the right (N->O, O->P, P->Q), leaving a blank column N sized to a width
# of 10 characters-worth of space.
#
# The active sheet/selection is also updated to match: "Repayment Schedule"
# becomes the active tab with its selection on R7 (previously "Transactions"
# was the active tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at N; existing N:P shift right to O:Q.
$ws.Columns("N:N").Insert() | Out-Null

# Size the freshly inserted column (renders as width 10 in the saved file).
$ws.Columns("N:N").ColumnWidth = 9.140625

# "Repayment Schedule" becomes the active sheet/tab, selection moves to R7.
# (Doing this last also naturally clears the previous tabSelected flag /
# selection state that belonged to the "Transactions" sheet.)
$ws.Activate() | Out-Null
$ws.Range("R7").Select() | Out-Null
